$wb = $excel.ActiveWorkbook

# =====================================================================
# Sheet "snapshot": refresh the injury list with the latest scrape.
# Two players recovered (removed) and two new players got injured
# (added); the list stays sorted by team_abbr, and every row gets a
# fresh "scraped_at" timestamp.
# =====================================================================
$ws1 = $wb.Worksheets.Item("snapshot")

$snapshotData = @(
    @("АВТ", "Автомобилист", "avtomobilist", "Зборовский Сергей", "2", "защитник", "20989", "1369_АВТ_зборовскийсергей", "injured_active", "https://www.khl.ru/clubs/avtomobilist/team/", "2025-12-02T03:01:25.459812+00:00"),
    @("АДМ", "Адмирал", "admiral", "Дарьин Александр", "41", "нападающий", "26282", "1369_АДМ_дарьиналександр", "injured_active", "https://www.khl.ru/clubs/admiral/team/", "2025-12-02T03:01:28.334910+00:00"),
    @("АДМ", "Адмирал", "admiral", "Педан Руслан", "53", "защитник", "25353", "1369_АДМ_педанруслан", "injured_active", "https://www.khl.ru/clubs/admiral/team/", "2025-12-02T03:01:28.334940+00:00"),
    @("АКБ", "Ак Барс", "ak_bars", "Яруллин Альберт", "33", "защитник", "16365", "1369_АКБ_яруллинальберт", "injured_active", "https://www.khl.ru/clubs/ak_bars/team/", "2025-12-02T03:01:31.208067+00:00"),
    @("АМР", "Амур", "amur", "Абросимов Роман", "94", "защитник", "17968", "1369_АМР_абросимовроман", "injured_active", "https://www.khl.ru/clubs/amur/team/", "2025-12-02T03:01:34.168708+00:00"),
    @("БАР", "Барыс", "barys", "Уотерспун Тайлер", "26", "защитник", "45769", "1369_БАР_уотерспунтайлер", "injured_active", "https://www.khl.ru/clubs/barys/team/", "2025-12-02T03:01:36.614119+00:00"),
    @("ДИН", "Динамо М", "dynamo_msk", "Готовец Кирилл", "41", "защитник", "16034", "1369_ДИН_готовецкирилл", "injured_active", "https://www.khl.ru/clubs/dynamo_msk/team/", "2025-12-02T03:01:40.372474+00:00"),
    @("ЛАД", "Лада", "lada", "Обидин Андрей", "58", "нападающий", "21310", "1369_ЛАД_обидинандрей", "injured_active", "https://www.khl.ru/clubs/lada/team/", "2025-12-02T03:01:45.429444+00:00"),
    @("ЛАД", "Лада", "lada", "Ожгихин Алексей", "43", "нападающий", "23021", "1369_ЛАД_ожгихиналексей", "injured_active", "https://www.khl.ru/clubs/lada/team/", "2025-12-02T03:01:45.429473+00:00"),
    @("ЛОК", "Локомотив", "lokomotiv", "Сергеев Андрей", "99", "защитник", "15416", "1369_ЛОК_сергеевандрей", "injured_active", "https://www.khl.ru/clubs/lokomotiv/team/", "2025-12-02T03:01:47.740760+00:00"),
    @("ММГ", "Металлург Мг", "metallurg_mg", "Козлов Андрей Е", "39", "нападающий", "40899", "1369_ММГ_козловандрейе", "injured_active", "https://www.khl.ru/clubs/metallurg_mg/team/", "2025-12-02T03:01:50.551854+00:00"),
    @("ММГ", "Металлург Мг", "metallurg_mg", "Сиряцкий Александр", "74", "защитник", "42458", "1369_ММГ_сиряцкийалександр", "injured_active", "https://www.khl.ru/clubs/metallurg_mg/team/", "2025-12-02T03:01:50.551883+00:00"),
    @("НХК", "Нефтехимик", "neftekhimik", "Хлыстов Никита", "7", "защитник", "17881", "1369_НХК_хлыстовникита", "injured_active", "https://www.khl.ru/clubs/neftekhimik/team/", "2025-12-02T03:01:53.347529+00:00"),
    @("СЕВ", "Северсталь", "severstal", "Фомин Макар", "77", "защитник", "42087", "1369_СЕВ_фоминмакар", "injured_active", "https://www.khl.ru/clubs/severstal/team/", "2025-12-02T03:01:55.765227+00:00"),
    @("СЕВ", "Северсталь", "severstal", "Цицюра Владислав", "10", "нападающий", "23840", "1369_СЕВ_цицюравладислав", "injured_active", "https://www.khl.ru/clubs/severstal/team/", "2025-12-02T03:01:55.765257+00:00"),
    @("СИБ", "Сибирь", "sibir", "Аланов Егор", "21", "защитник", "26698", "1369_СИБ_алановегор", "injured_active", "https://www.khl.ru/clubs/sibir/team/", "2025-12-02T03:01:58.086071+00:00"),
    @("СКА", "СКА", "ska", "Зайцев Никита И", "22", "защитник", "16024", "1369_СКА_зайцевникитаи", "injured_active", "https://www.khl.ru/clubs/ska/team/", "2025-12-02T03:02:00.847732+00:00"),
    @("СКА", "СКА", "ska", "Короткий Матвей", "71", "нападающий", "41566", "1369_СКА_короткийматвей", "injured_active", "https://www.khl.ru/clubs/ska/team/", "2025-12-02T03:02:00.847763+00:00"),
    @("СОЧ", "ХК Сочи", "hc_sochi", "Бикмуллин Рафаэль", "24", "нападающий", "22424", "1369_СОЧ_бикмуллинрафаэль", "injured_active", "https://www.khl.ru/clubs/hc_sochi/team/", "2025-12-02T03:02:03.638482+00:00"),
    @("СОЧ", "ХК Сочи", "hc_sochi", "Венгрыжановский Денис", "9", "нападающий", "31892", "1369_СОЧ_венгрыжановскийденис", "injured_active", "https://www.khl.ru/clubs/hc_sochi/team/", "2025-12-02T03:02:03.638513+00:00"),
    @("СОЧ", "ХК Сочи", "hc_sochi", "Хёфенмайер Ноэль", "22", "защитник", "44847", "1369_СОЧ_хефенмайерноэль", "injured_active", "https://www.khl.ru/clubs/hc_sochi/team/", "2025-12-02T03:02:03.638532+00:00"),
    @("СПР", "Спартак", "spartak", "Вишневский Дмитрий", "55", "защитник", "15299", "1369_СПР_вишневскийдмитрий", "injured_active", "https://www.khl.ru/clubs/spartak/team/", "2025-12-02T03:02:06.436208+00:00"),
    @("ТОР", "Торпедо", "torpedo", "Кручинин Алексей", "78", "нападающий", "16355", "1369_ТОР_кручининалексей", "injured_active", "https://www.khl.ru/clubs/torpedo/team/", "2025-12-02T03:02:11.669619+00:00"),
    @("ТРК", "Трактор", "traktor", "Мыльников Сергей И", "20", "вратарь", "24799", "1369_ТРК_мыльниковсергейи", "injured_active", "https://www.khl.ru/clubs/traktor/team/", "2025-12-02T03:02:14.398134+00:00"),
    @("ТРК", "Трактор", "traktor", "Светлаков Андрей", "87", "нападающий", "19218", "1369_ТРК_светлаковандрей", "injured_active", "https://www.khl.ru/clubs/traktor/team/", "2025-12-02T03:02:14.398163+00:00"),
    @("ЦСК", "ЦСКА", "cska", "Бучельников Дмитрий", "72", "нападающий", "39102", "1369_ЦСК_бучельниковдмитрий", "injured_active", "https://www.khl.ru/clubs/cska/team/", "2025-12-02T03:02:17.317732+00:00"),
    @("ЦСК", "ЦСКА", "cska", "Моисеев Данила", "93", "нападающий", "23931", "1369_ЦСК_моисеевданила", "injured_active", "https://www.khl.ru/clubs/cska/team/", "2025-12-02T03:02:17.317832+00:00"),
    @("ШДР", "Драконы", "kunlun", "Бишофф Джейк", "28", "защитник", "45490", "1369_ШДР_бишоффджейк", "injured_active", "https://www.khl.ru/clubs/kunlun/team/", "2025-12-02T03:02:19.602937+00:00"),
    @("ШДР", "Драконы", "kunlun", "Гроло Жереми", "75", "защитник", "45343", "1369_ШДР_гроложереми", "injured_active", "https://www.khl.ru/clubs/kunlun/team/", "2025-12-02T03:02:19.602965+00:00"),
)

$snapshotRowCount = $snapshotData.Count
$snapshotColCount = $snapshotData[0].Count

# Columns E (number) and G (player_id_khl) hold digit-only text such as
# "2" or "20989"; force them to Text format so Excel keeps them as
# strings instead of silently converting them to numbers.
$numTxtRange = $ws1.Range("E2:E" + (1 + $snapshotRowCount))
$numTxtRange.NumberFormat = "@"
$idTxtRange = $ws1.Range("G2:G" + (1 + $snapshotRowCount))
$idTxtRange.NumberFormat = "@"

for ($i = 0; $i -lt $snapshotRowCount; $i++) {
    $r = $i + 2
    $rowVals = $snapshotData[$i]
    for ($c = 0; $c -lt $snapshotColCount; $c++) {
        $ws1.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

# =====================================================================
# Sheet "returned": the two players who recovered since the last run.
# The "status" column (E) already reads "RETURN" and is left as-is;
# we just refresh who they are and when the change was detected.
# =====================================================================
$ws2 = $wb.Worksheets.Item("returned")

$returnedData = @(
    @("СЕВ", "Северсталь", "Квочко Илья", "1369_СЕВ_квочкоилья", "RETURN", "2025-12-02T11:02:20.114773+08:00", "2025-12-02"),
    @("СИБ", "Сибирь", "Пьянов Валентин", "1369_СИБ_пьяноввалентин", "RETURN", "2025-12-02T11:02:20.114773+08:00", "2025-12-02"),
)

# Column G (changed_day) looks like a date ("2025-12-02"); keep it text.
$ws2.Range("G2:G" + (1 + $returnedData.Count)).NumberFormat = "@"

for ($i = 0; $i -lt $returnedData.Count; $i++) {
    $r = $i + 2
    $rowVals = $returnedData[$i]
    # Column E (status) is intentionally left untouched (stays "RETURN").
    $ws2.Cells.Item($r, 1).Value = $rowVals[0]
    $ws2.Cells.Item($r, 2).Value = $rowVals[1]
    $ws2.Cells.Item($r, 3).Value = $rowVals[2]
    $ws2.Cells.Item($r, 4).Value = $rowVals[3]
    $ws2.Cells.Item($r, 6).Value = $rowVals[5]
    $ws2.Cells.Item($r, 7).Value = $rowVals[6]
}

# =====================================================================
# Sheet "new_injured": the two players newly added to the injury list.
# The sheet previously only had the header row.
# =====================================================================
$ws3 = $wb.Worksheets.Item("new_injured")

$newInjuredData = @(
    @("ТОР", "Торпедо", "Кручинин Алексей", "1369_ТОР_кручининалексей", "INJURED_NEW", "2025-12-02T11:02:20.114773+08:00", "2025-12-02"),
    @("ТРК", "Трактор", "Светлаков Андрей", "1369_ТРК_светлаковандрей", "INJURED_NEW", "2025-12-02T11:02:20.114773+08:00", "2025-12-02"),
)

# Column G (changed_day) looks like a date ("2025-12-02"); keep it text.
$ws3.Range("G2:G" + (1 + $newInjuredData.Count)).NumberFormat = "@"

for ($i = 0; $i -lt $newInjuredData.Count; $i++) {
    $r = $i + 2
    $rowVals = $newInjuredData[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws3.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
